{"js": "// Remove the four paragraphs that follow the \"LOB1036: ...\" requirement\n// paragraph: an empty paragraph, the \"Ver no Jupiter...\" paragraph, another\n// empty paragraph, and the empty page-break paragraph that had jc=\"left\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet targetIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"LOB1036: Geometria Anal\u00edtica (Requisito fraco)\") {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error(\"Could not find the 'LOB1036' requirement paragraph\");\n}\n\n// Delete the next four paragraphs (indices targetIndex+1 .. targetIndex+4)\nconst toDelete = [];\nfor (let i = targetIndex + 1; i <= targetIndex + 4 && i < items.length; i++) {\n  toDelete.push(items[i]);\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the four paragraphs that follow the \"LOB1036: ...\" requirement\n# paragraph: an empty paragraph, the \"Ver no Jupiter...\" paragraph, another\n# empty paragraph, and the empty page-break paragraph that had jc=\"left\".\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"LOB1036:*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'LOB1036' requirement paragraph\"\n}\n\nfor ($k = 0; $k -lt 4; $k++) {\n    $next = $target.Next()\n    $next.Range.Delete()\n}\n"}
